# Update the "Keyword" (column B) and "Correlation" (column C) values on the
# active worksheet to reflect the corrected/shuffled correlation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "credit"
$ws.Range("C2").Value = 0.1327

$ws.Range("B3").Value = "trade"
$ws.Range("C3").Value = 0.2366

$ws.Range("B4").Value = "bank"
$ws.Range("C4").Value = -0.1438

$ws.Range("B7").Value = "bank"
$ws.Range("C7").Value = -0.2167

$ws.Range("B8").Value = "credit"
$ws.Range("C8").Value = -0.5286999999999999

$ws.Range("B9").Value = "inflation"
$ws.Range("C9").Value = -0.1311

$ws.Range("B10").Value = "trade"
$ws.Range("C10").Value = 0.1349

$ws.Range("B12").Value = "credit"
$ws.Range("C12").Value = -0.4083

$ws.Range("B13").Value = "trade"
$ws.Range("C13").Value = 0.8954

$ws.Range("B14").Value = "bank"
$ws.Range("C14").Value = -0.2377

$ws.Range("B17").Value = "bank"
$ws.Range("C17").Value = 0.1918

$ws.Range("B18").Value = "trade"
$ws.Range("C18").Value = -0.1808

$ws.Range("B19").Value = "inflation"
$ws.Range("C19").Value = 0.5832000000000001
